$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.09709999999999
$ws.Range("B4").Value = 8.544899999999995

$ws.Range("A7").Value = -20.26029999999997

$ws.Range("D10").Value = -7.428200000000004

$ws.Range("B12").Value = 5.329799999999998

$ws.Range("D13").Value = -8.320699999999993

$ws.Range("A16").Value = -22.01480000000002

$ws.Range("B18").Value = 6.677199999999999

$ws.Range("B19").Value = 8.568100000000005

$ws.Range("B20").Value = 9.57279999999999

$ws.Range("A28").Value = -21.8789

$ws.Range("A29").Value = -21.18039999999997

$ws.Range("D30").Value = -7.026699999999995

$ws.Range("B31").Value = 4.607799999999998

$ws.Range("A32").Value = -21.1748

$ws.Range("A40").Value = -20.28509999999999
$ws.Range("B40").Value = 8.414600000000004
$ws.Range("D40").Value = -8.041000000000006

$ws.Range("B42").Value = 8.649299999999998

$ws.Range("D44").Value = -6.839300000000001

$ws.Range("B47").Value = 4.949900000000002

$ws.Range("B48").Value = 7.261800000000006

$ws.Range("A52").Value = -22.2507

$ws.Range("A57").Value = -22.00690000000002

$ws.Range("B63").Value = 4.865499999999997

$ws.Range("B64").Value = 5.943700000000002

$ws.Range("A66").Value = -21.4149

$ws.Range("B76").Value = 5.782899999999996

$ws.Range("B81").Value = 5.017100000000004

$ws.Range("B89").Value = 4.732599999999992
$ws.Range("D89").Value = -8.544199999999998

$ws.Range("D91").Value = -7.6754

$ws.Range("B94").Value = 4.884099999999993

$ws.Range("A100").Value = -22.04320000000003
